# Generate Report for Handback
# Refreshes the timestamp strings recorded on the handback status report.
#
#   Overview!G2  "Latest HO Xliff Generate Date"        for 607ed64b-...md
#   zh-cn!H2     "Correspond Handoff Datetime"           for 607ed64b-...md
#   zh-cn!K2     "Correspond Handback DateTime"          for 607ed64b-...md
#   de-de!H2     "Correspond Handoff Datetime"           for 607ed64b-...md
#   de-de!K2     "Correspond Handback DateTime"          for 607ed64b-...md

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet - Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-09-07 03:22:23"

# zh-cn sheet - Correspond Handoff / Handback Datetime
$wsZhCn.Range("H2").Value = "2016-09-07 03:22:17"
$wsZhCn.Range("K2").Value = "2016-09-07 03:22:35"

# de-de sheet - Correspond Handoff / Handback Datetime
$wsDeDe.Range("H2").Value = "2016-09-07 03:22:23"
$wsDeDe.Range("K2").Value = "2016-09-07 03:22:43"
